# Weekly price-sheet refresh: a new price observation is inserted as row 38
# ("Berenjena", Primera, 2023-05-03) and every existing record from the old
# row 38 onward shifts down by one row (old row 72 becomes the new row 73).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 38, pushing rows 38:72 down to 39:73.
$ws.Rows(38).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(38, 1).Value  = 7
$ws.Cells.Item(38, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38, 3).Value  = "Ñuble"
$ws.Cells.Item(38, 4).Value  = 45049
$ws.Cells.Item(38, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(38, 5).Value  = 16
$ws.Cells.Item(38, 6).Value  = 100112001
$ws.Cells.Item(38, 7).Value  = "Berenjena"
$ws.Cells.Item(38, 8).Value  = "Sin especificar"
$ws.Cells.Item(38, 9).Value  = "Primera"
$ws.Cells.Item(38, 10).Value = 70
$ws.Cells.Item(38, 11).Value = 10000
$ws.Cells.Item(38, 12).Value = 11000
$ws.Cells.Item(38, 13).Value = 10571
$ws.Cells.Item(38, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(38, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(38, 16).Value = 176
$ws.Cells.Item(38, 17).Value = 60
$ws.Cells.Item(38, 18).Value = "Hortaliza"
